$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1730205278592375
$ws.Range("C2").Value = 0.5894428152492669
$ws.Range("J2").Value = 0.008797653958944282
$ws.Range("P2").Value = 0.1378299120234604
$ws.Range("S2").Value = 0.09090909090909091
$ws.Range("C3").Value = 0.01463414634146342
$ws.Range("J3").Value = 0.02439024390243903
$ws.Range("P3").Value = 0.7121951219512195
$ws.Range("S3").Value = 0.248780487804878
$ws.Range("J4").Value = 0.075
$ws.Range("P4").Value = 0.7
$ws.Range("S4").Value = 0.225
$ws.Range("B6").Value = 0.05860805860805861
$ws.Range("D6").Value = 0.01465201465201465
$ws.Range("E6").Value = 0.003663003663003663
$ws.Range("F6").Value = 0.06227106227106227
$ws.Range("J6").Value = 0.3003663003663004
$ws.Range("O6").Value = 0.01831501831501832
$ws.Range("Q6").Value = 0.1025641025641026
$ws.Range("R6").Value = 0.0695970695970696
$ws.Range("S6").Value = 0.36996336996337
$ws.Range("B7").Value = 0.08994708994708994
$ws.Range("D7").Value = 0.01058201058201058
$ws.Range("F7").Value = 0.04761904761904762
$ws.Range("J7").Value = 0.1322751322751323
$ws.Range("O7").Value = 0.02116402116402116
$ws.Range("Q7").Value = 0.1481481481481481
$ws.Range("R7").Value = 0.08994708994708994
$ws.Range("S7").Value = 0.4603174603174603
$ws.Range("B8").Value = 0.1146067415730337
$ws.Range("D8").Value = 0.01123595505617977
$ws.Range("F8").Value = 0.0651685393258427
$ws.Range("J8").Value = 0.0898876404494382
$ws.Range("O8").Value = 0.02921348314606742
$ws.Range("Q8").Value = 0.1528089887640449
$ws.Range("R8").Value = 0.07865168539325842
$ws.Range("S8").Value = 0.4584269662921348
$ws.Range("B9").Value = 0.1061946902654867
$ws.Range("D9").Value = 0.004424778761061947
$ws.Range("E9").Value = 0.004424778761061947
$ws.Range("F9").Value = 0.09292035398230089
$ws.Range("J9").Value = 0.1106194690265487
$ws.Range("O9").Value = 0.008849557522123894
$ws.Range("Q9").Value = 0.1504424778761062
$ws.Range("R9").Value = 0.0752212389380531
$ws.Range("S9").Value = 0.4469026548672566
$ws.Range("B10").Value = 0.1371610845295056
$ws.Range("D10").Value = 0.02312599681020734
$ws.Range("E10").Value = 0.002392344497607655
$ws.Range("F10").Value = 0.08133971291866028
$ws.Range("J10").Value = 0.09569377990430622
$ws.Range("O10").Value = 0.01913875598086124
$ws.Range("Q10").Value = 0.1690590111642743
$ws.Range("R10").Value = 0.08293460925039872
$ws.Range("S10").Value = 0.3891547049441786
$ws.Range("G11").Value = 0.1757188498402556
$ws.Range("J11").Value = 0.07987220447284345
$ws.Range("K11").Value = 0.2364217252396166
$ws.Range("L11").Value = 0.4984025559105431
$ws.Range("S11").Value = 0.009584664536741214
$ws.Range("G12").Value = 0.717948717948718
$ws.Range("J12").Value = 0.217948717948718
$ws.Range("K12").Value = 0.01282051282051282
$ws.Range("L12").Value = 0.01923076923076923
$ws.Range("S12").Value = 0.03205128205128205
$ws.Range("G13").Value = 0.6304347826086957
$ws.Range("J13").Value = 0.3260869565217391
$ws.Range("S13").Value = 0.04347826086956522
$ws.Range("F15").Value = 0.02531645569620253
$ws.Range("H15").Value = 0.1308016877637131
$ws.Range("I15").Value = 0.09282700421940929
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.04641350210970464
$ws.Range("M15").Value = 0.02109704641350211
$ws.Range("O15").Value = 0.1012658227848101
$ws.Range("S15").Value = 0.2489451476793249
$ws.Range("F16").Value = 0.01363636363636364
$ws.Range("H16").Value = 0.1772727272727273
$ws.Range("I16").Value = 0.1045454545454545
$ws.Range("J16").Value = 0.3318181818181818
$ws.Range("K16").Value = 0.1272727272727273
$ws.Range("M16").Value = 0.03181818181818181
$ws.Range("O16").Value = 0.08636363636363636
$ws.Range("S16").Value = 0.1272727272727273
$ws.Range("F17").Value = 0.02419354838709677
$ws.Range("H17").Value = 0.2043010752688172
$ws.Range("I17").Value = 0.08870967741935484
$ws.Range("J17").Value = 0.4301075268817204
$ws.Range("K17").Value = 0.08333333333333333
$ws.Range("M17").Value = 0.02150537634408602
$ws.Range("O17").Value = 0.06182795698924731
$ws.Range("S17").Value = 0.08602150537634409
$ws.Range("F18").Value = 0.03125
$ws.Range("H18").Value = 0.15625
$ws.Range("I18").Value = 0.09375
$ws.Range("J18").Value = 0.4739583333333333
$ws.Range("K18").Value = 0.078125
$ws.Range("M18").Value = 0.015625
$ws.Range("O18").Value = 0.07291666666666667
$ws.Range("S18").Value = 0.078125
$ws.Range("F19").Value = 0.03305785123966942
$ws.Range("H19").Value = 0.2096168294515402
$ws.Range("I19").Value = 0.09691960931630353
$ws.Range("J19").Value = 0.3703981968444778
$ws.Range("K19").Value = 0.1104432757325319
$ws.Range("M19").Value = 0.01953418482344102
$ws.Range("O19").Value = 0.06160781367392937
$ws.Range("S19").Value = 0.09842223891810668
